# Case_7_27 diagnostic.xlsx edit:
#  - B1 = 0, A2 = 0, both styled bold/centered/top-aligned with a thin box border
#  - B2 = "disconnected_elements" (plain, unstyled)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the shared format (bold font, thin box border, centered/top aligned)
# on B1 first ...
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108  # xlCenter
$r1.VerticalAlignment = -4160    # xlTop
$r1.Borders.LineStyle = 1        # xlContinuous
$r1.Borders.Weight = 2           # xlThin

# ... then copy/paste the resulting format onto A2 so both cells land on the
# exact same style record (setting the properties a second time from scratch
# would otherwise mint a duplicate, transient style entry).
$r1.Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
